$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after the existing one, named
# "description_treatment_arms" (note: no trailing "s", unlike the
# original sheet "description_treatment_armss").
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "description_treatment_arms"

# Populate it with the numeric (# cases / # treated cases) and date
# columns copied from the first sheet, using a value-only paste so the
# shared-string / numeric cell types (and lack of any special number
# formatting) are preserved exactly as on sheet1 rather than being
# re-inferred (which would turn the dd/mm/yyyy text into real dates).
$ws1.Range("D2:E4").Copy()
$ws2.Range("D2:E4").PasteSpecial(-4163)

$ws1.Range("K2:L4").Copy()
$ws2.Range("K2:L4").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# Keep the original sheet as the active / selected one.
$ws1.Activate()
[void]$ws1.Range("A1:H4").Select()

$wb.ForceFullCalculation = $true
$excel.CalculateFullRebuild()
